$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two lanes ("faixas") that had been left out of the sheet: insert two
# rows at row 28 (pushing the existing rows 28-32 down to rows 30-34).
$ws.Rows.Item(28).Resize(2).Insert()

# New row 28: AVENIDA BRAZ LEME
$ws.Range("A28").Value = "AVENIDA BRAZ LEME"
$ws.Range("B28").Value = "AVENIDA BRAZ LEME"
$ws.Range("C28").Value2 = 29
$ws.Range("D28").Value2 = 5
$ws.Range("E28").Value2 = 2024
$ws.Range("I28").Value2 = 23
$ws.Range("J28").Value2 = 11

# New row 29: AVENIDA ELISEU DE ALMEIDA
$ws.Range("A29").Value = "AVENIDA ELISEU DE ALMEIDA"
$ws.Range("B29").Value = "AVENIDA ELISEU DE ALMEIDA"
$ws.Range("C29").Value2 = 29
$ws.Range("D29").Value2 = 5
$ws.Range("E29").Value2 = 2024
$ws.Range("I29").Value2 = 24
$ws.Range("J29").Value2 = 12

# The rows that got shifted down (old 28-32, now 30-34) keep their data,
# but id_logradouro (I) and id_trecho (J) need to be renumbered to follow
# on from the two newly inserted rows.
$ws.Range("I30").Value2 = 25
$ws.Range("I31").Value2 = 26
$ws.Range("I32").Value2 = 27
$ws.Range("I33").Value2 = 28
$ws.Range("I34").Value2 = 29

$ws.Range("J32").Value2 = 13
$ws.Range("J33").Value2 = 14
$ws.Range("J34").Value2 = 15

# Restore the view state (frozen-pane top-left cell and active selection).
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("B20").Select()
